# This script applies the "gh-pages output" refresh for 合肥-漫展信息.xlsx:
#  - Removes the cancelled event row (合肥·WA二次元饭局（取消）, 2024-05-18) from the
#    "展览" (sheet 1) and "全部类型" (sheet 4) worksheets, shifting all following
#    rows up by one.
#  - Re-numbers the index column (A) after the shift.
#  - Refreshes the "想去人数" (F) / "最低票价" (G) counters that bilibili reported
#    for the remaining events.

$wb = $excel.ActiveWorkbook

function Update-ComicSheet {
    param($Worksheet, $Overrides)

    # Drop the cancelled "WA二次元饭局" row; Excel shifts every following row up.
    $Worksheet.Rows(2).Delete() | Out-Null

    # Column A holds hard-coded sequence numbers (1, 2, 3, ...) rather than a
    # formula, so it needs to be re-sequenced after the shift.
    $lastRow = $Worksheet.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $Worksheet.Cells.Item($r, 1).Value = $r - 1
    }

    # Refresh the scraped "want to go" counts / lowest ticket price for the
    # events whose numbers changed since the last crawl.
    foreach ($rowKey in $Overrides.Keys) {
        $rowOverrides = $Overrides[$rowKey]
        foreach ($colKey in $rowOverrides.Keys) {
            $Worksheet.Cells.Item([int]$rowKey, [int]$colKey).Value = $rowOverrides[$colKey]
        }
    }
}

# Column map: F = 6 (想去人数), G = 7 (最低票价)
$overrides1 = @{
    3  = @{ 6 = 45 }
    4  = @{ 6 = 606 }
    6  = @{ 6 = 9105 }
    7  = @{ 6 = 829 }
    9  = @{ 6 = 1168; 7 = 29.9 }
    10 = @{ 6 = 1048 }
    11 = @{ 6 = 132 }
    12 = @{ 6 = 51 }
    13 = @{ 6 = 12 }
    14 = @{ 6 = 245 }
    15 = @{ 6 = 339 }
    17 = @{ 6 = 242 }
    18 = @{ 6 = 1175 }
}

$overrides4 = @{
    4  = @{ 6 = 45 }
    6  = @{ 6 = 606 }
    8  = @{ 6 = 9105 }
    9  = @{ 6 = 829 }
    11 = @{ 6 = 1168; 7 = 29.9 }
    12 = @{ 6 = 1048 }
    13 = @{ 6 = 132 }
    14 = @{ 6 = 51 }
    15 = @{ 6 = 12 }
    16 = @{ 6 = 245 }
    17 = @{ 6 = 339 }
    19 = @{ 6 = 242 }
    20 = @{ 6 = 1175 }
}

# Sheet 1: 展览 (Exhibitions)
Update-ComicSheet $wb.Worksheets.Item(1) $overrides1

# Sheet 4: 全部类型 (All types)
Update-ComicSheet $wb.Worksheets.Item(4) $overrides4
